$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, whether the value must be
# force-written as Text (to stop Excel from auto-converting a numeric-
# looking string like "0.999" or "142.45" into a real number).
$updates = @(
    ,@("D2", "68.668.10", $false)
    ,@("E2", "  +1.11%  ", $false)
    ,@("D3", "3.872.99", $false)
    ,@("E3", "  +0.57%  ", $false)
    ,@("D4", "0.999", $true)
    ,@("E4", "  -0.05%  ", $false)
    ,@("D5", "602.88", $true)
    ,@("E5", "  +0.70%  ", $false)
    ,@("D6", "171.92", $true)
    ,@("E6", "  +3.37%  ", $false)
    ,@("D7", "3.873.35", $false)
    ,@("E7", "  +0.59%  ", $false)
    ,@("E8", "  +0.03%  ", $false)
    ,@("E9", "  +0.98%  ", $false)
    ,@("E10", "  +2.53%  ", $false)
    ,@("D11", "6.53", $true)
    ,@("E11", "  +3.65%  ", $false)
    ,@("E12", "  +1.50%  ", $false)
    ,@("D13", "0.0000288", $true)
    ,@("E13", "  +16.07%  ", $false)
    ,@("D14", "37.24", $true)
    ,@("E14", "  +0.96%  ", $false)
    ,@("D15", "4.520.70", $false)
    ,@("E15", "  +0.53%  ", $false)
    ,@("D16", "3.879.83", $false)
    ,@("E16", "  +0.57%  ", $false)
    ,@("D17", "68.766.97", $false)
    ,@("E17", "  +1.22%  ", $false)
    ,@("D18", "18.33", $true)
    ,@("E18", "  +0.68%  ", $false)
    ,@("E19", "  +0.26%  ", $false)
    ,@("E20", "  +1.07%  ", $false)
    ,@("D21", "10.99", $true)
    ,@("E21", "  +2.01%  ", $false)
    ,@("D22", "472.82", $true)
    ,@("E22", "  +1.29%  ", $false)
    ,@("D23", "0.734", $true)
    ,@("E23", "  +0.81%  ", $false)
    ,@("D24", "0.0000163", $true)
    ,@("E24", "  +1.88%  ", $false)
    ,@("D25", "84.04", $true)
    ,@("E25", "  +1.01%  ", $false)
    ,@("E26", "  +2.55%  ", $false)
    ,@("D27", "12.28", $true)
    ,@("E27", "  +1.10%  ", $false)
    ,@("D28", "10.50", $true)
    ,@("E28", "  +4.94%  ", $false)
    ,@("E29", "  +0.09%  ", $false)
    ,@("E30", "  +0.41%  ", $false)
    ,@("D31", "4.024.00", $false)
    ,@("E31", "  +0.57%  ", $false)
    ,@("D32", "7.78", $true)
    ,@("E32", "  +1.36%  ", $false)
    ,@("E33", "  +1.13%  ", $false)
    ,@("D34", "31.38", $true)
    ,@("E34", "  +1.50%  ", $false)
    ,@("E35", "  +0.48%  ", $false)
    ,@("D36", "3.838.15", $false)
    ,@("E37", "  +22.09%  ", $false)
    ,@("E38", "  +1.41%  ", $false)
    ,@("E39", "  +2.25%  ", $false)
    ,@("D40", "0.141", $true)
    ,@("E40", "  +0.65%  ", $false)
    ,@("E41", "  +1.05%  ", $false)
    ,@("E42", "  +0.09%  ", $false)
    ,@("D43", "0.322", $true)
    ,@("E43", "  +3.07%  ", $false)
    ,@("D44", "0.000304", $true)
    ,@("E44", "  +12.05%  ", $false)
    ,@("E45", "  +1.17%  ", $false)
    ,@("E46", "  -0.02%  ", $false)
    ,@("E47", "  +2.95%  ", $false)
    ,@("D48", "420.45", $true)
    ,@("E48", "  -1.44%  ", $false)
    ,@("D49", "46.70", $true)
    ,@("E49", "  -0.87%  ", $false)
    ,@("B50", "VeChain", $false)
    ,@("C50", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", $false)
    ,@("D50", "0.0361", $true)
    ,@("E50", "  +1.97%  ", $false)
    ,@("B51", "Monero", $false)
    ,@("C51", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", $false)
    ,@("D51", "142.45", $true)
    ,@("E51", "  -0.68%  ", $false)
)

foreach ($item in $updates) {
    $addr = $item[0]
    $val = $item[1]
    $forceText = $item[2]
    $rng = $ws.Range($addr)
    if ($forceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $val
}
